# Applies crypto price/volume updates per the commit diff (Wed Dec 13 09:37:14 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.141.59'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '2.176.01'
$ws.Range('E3').Value = '  -2.40%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '249.61'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -2.50%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '66.55'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -7.44%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.575'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '59.08'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '36.55'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -11.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0936'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.45%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.103'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.88'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.20%  '
$ws.Range('D15').Value = '2.503.71'
$ws.Range('E15').Value = '  -2.33%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.28'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.74%  '
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '2.182.55'
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('D19').Value = '41.109.82'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').Value = '0.0₃0951'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.77'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '229.05'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('E25').Value = '  -7.42%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.45'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.82%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.40'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.14%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.72'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.84%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.30'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('E31').Value = '  -7.08%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.25'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.122'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.79'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0763'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.98%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.122'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.11%  '
$ws.Range('E37').Value = '  -4.28%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.99'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '24.54'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -6.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0308'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.74%  '
$ws.Range('E41').Value = '  -3.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.29'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.48%  '
$ws.Range('E43').Value = '  -8.39%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '61.24'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -9.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '11.37'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.29%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.52'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  -6.98%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.14'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.21%  '
